# chore: adapt column header formatting to respective input file names
#
# 1. Rename the two header blocks (cols A-J and L-U) from the generic
#    "_old" / "_new" suffixes to the format-version-specific
#    "_FV2210" / "_FV2304" suffixes (column K / "diff" stays as-is).
# 2. Turn the data range into a real Excel Table ("Table1").
# 3. Freeze the header row (pane split below row 1).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Header renames ------------------------------------------------

$fv2210Headers = @(
    "Segmentname_FV2210",
    "Segmentgruppe_FV2210",
    "Segment_FV2210",
    "Datenelement_FV2210",
    "Segment ID_FV2210",
    "Code_FV2210",
    "Qualifier_FV2210",
    "Beschreibung_FV2210",
    "Bedingungsausdruck_FV2210",
    "Bedingung_FV2210"
)

for ($i = 0; $i -lt $fv2210Headers.Length; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $fv2210Headers[$i]
}

# Column K (11) holds "diff" and is left untouched.

$fv2304Headers = @(
    "Segmentname_FV2304",
    "Segmentgruppe_FV2304",
    "Segment_FV2304",
    "Datenelement_FV2304",
    "Segment ID_FV2304",
    "Code_FV2304",
    "Qualifier_FV2304",
    "Beschreibung_FV2304",
    "Bedingungsausdruck_FV2304",
    "Bedingung_FV2304"
)

for ($i = 0; $i -lt $fv2304Headers.Length; $i++) {
    $ws.Cells.Item(1, $i + 12).Value = $fv2304Headers[$i]
}

# --- 2. Convert the used range into an Excel Table ---------------------

$dataRange = $ws.Range("A1:U67")
$tbl = $ws.ListObjects.Add(
    [Microsoft.Office.Interop.Excel.XlListObjectSourceType]::xlSrcRange,
    $dataRange,
    $null,
    [Microsoft.Office.Interop.Excel.XlYesNoGuess]::xlYes
)
$tbl.Name = "Table1"

# --- 3. Freeze header row ------------------------------------------------

$ws.Activate()
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true

Write-Host "Header rename, table conversion and freeze panes applied."
